$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 13; this shifts all existing rows 13:98 down to 14:99
# and extends the used range to row 99.
$ws.Rows.Item(13).Insert()

# Populate the newly inserted row 13 with the new weekly record.
$ws.Cells.Item(13, 1).Value = 1
$ws.Cells.Item(13, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(13, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(13, 4).Value = 44687
$ws.Cells.Item(13, 5).Value = 15
$ws.Cells.Item(13, 6).Value = "Fruta"
$ws.Cells.Item(13, 7).Value = 100102
$ws.Cells.Item(13, 8).Value = "Cítricos"
$ws.Cells.Item(13, 9).Value = 100102004
$ws.Cells.Item(13, 10).Value = "Mandarina"
$ws.Cells.Item(13, 11).Value = "Murcott"
$ws.Cells.Item(13, 12).Value = "Primera"
$ws.Cells.Item(13, 13).Value = 300
$ws.Cells.Item(13, 14).Value = 19000
$ws.Cells.Item(13, 15).Value = 20000
$ws.Cells.Item(13, 16).Value = 19500
$ws.Cells.Item(13, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(13, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(13, 19).Value = 975
$ws.Cells.Item(13, 20).Value = 20
